# "output for led driver correct on scope for cmsis"
#
# Sheet1: B1's driver-frequency formula becomes what used to be in C1
# (2.7*10^6 instead of 2.6*10^6); C1 and D1 (the other two candidate
# frequencies) are no longer needed and are cleared. Everything
# downstream (B2, I1:I5, E/G columns for rows 7-13) recalculates off of
# the new B1 automatically. The scratch scope readings in H9:J9 and the
# stale H10:I10 readings are removed since they're superseded by the
# corrected output.
#
# View state: Sheet1 becomes the active/selected tab (zoomed to 145%,
# cursor on J10) instead of Sheet2 (cursor left on G15).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 data edits ---------------------------------------------------

# B1 now carries the 2.7MHz candidate that used to live in C1.
$ws1.Range("B1").Formula = "=2.7*10^6"

# C1 / D1 (2.6MHz / 2.4MHz candidates) are no longer used.
$ws1.Range("C1").ClearContents()
$ws1.Range("D1").ClearContents()

# Stale scope-measurement values, superseded by the corrected output.
$ws1.Range("H9:J9").ClearContents()
$ws1.Range("H10:I10").ClearContents()

# --- View / selection state ----------------------------------------------

# Leave Sheet2's selection on G15 without leaving it the active tab:
# activate it first, move the selection, then activate Sheet1 last so
# Sheet1 ends up as the active/selected sheet.
$ws2.Activate()
$ws2.Range("G15").Select()

$ws1.Activate()
$excel.ActiveWindow.Zoom = 145
$ws1.Range("J10").Select()
